$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "65.851.99"
$ws.Range("D3").Value = "2.694.68"
$ws.Range("E3").Value = "  +2.12%  "
$ws.Range("E4").Value = "  -0.01%  "
$ws.Range("D5").Value = "'607.73"
$ws.Range("D5").Style = $ws.Range("B5").Style
$ws.Range("D6").Value = "'157.75"
$ws.Range("D6").Style = $ws.Range("B6").Style
$ws.Range("E6").Value = "  +1.49%  "
$ws.Range("E7").Value = "  -0.01%  "
$ws.Range("D8").Value = "'0.588"
$ws.Range("D8").Style = $ws.Range("B8").Style
$ws.Range("E8").Value = "  -0.73%  "
$ws.Range("E9").Value = "  +5.79%  "
$ws.Range("D10").Value = "'6.01"
$ws.Range("D10").Style = $ws.Range("B10").Style
$ws.Range("E10").Value = "  +3.70%  "
$ws.Range("E11").Value = "  +0.15%  "
$ws.Range("E12").Value = "  +1.11%  "
$ws.Range("D13").Value = "'30.32"
$ws.Range("D13").Style = $ws.Range("B13").Style
$ws.Range("D14").Value = "'0.0000202"
$ws.Range("D14").Style = $ws.Range("B14").Style
$ws.Range("E14").Value = "  +7.56%  "
$ws.Range("D15").Value = "3.177.97"
$ws.Range("E15").Value = "  +2.07%  "
$ws.Range("D16").Value = "65.681.82"
$ws.Range("E16").Value = "  +1.19%  "
$ws.Range("D17").Value = "2.688.81"
$ws.Range("E17").Value = "  +0.43%  "
$ws.Range("E18").Value = "  +0.48%  "
$ws.Range("D19").Value = "'4.89"
$ws.Range("D19").Style = $ws.Range("B19").Style
$ws.Range("E19").Value = "  +2.04%  "
$ws.Range("D20").Value = "'358.71"
$ws.Range("D20").Style = $ws.Range("B20").Style
$ws.Range("E20").Value = "  +1.79%  "
$ws.Range("E21").Value = "  +2.96%  "
$ws.Range("E22").Value = "  -0.70%  "
$ws.Range("D23").Value = "'70.59"
$ws.Range("D23").Style = $ws.Range("B23").Style
$ws.Range("E23").Value = "  +3.74%  "
$ws.Range("D24").Value = "'9.84"
$ws.Range("D24").Style = $ws.Range("B24").Style
$ws.Range("E24").Value = "  +2.91%  "
$ws.Range("D25").Value = "'1.68"
$ws.Range("D25").Style = $ws.Range("B25").Style
$ws.Range("E25").Value = "  -0.82%  "
$ws.Range("D26").Value = "'0.0000107"
$ws.Range("D26").Style = $ws.Range("B26").Style
$ws.Range("E26").Value = "  +13.26%  "
$ws.Range("E27").Value = "  +2.49%  "
$ws.Range("E28").Value = "  +4.97%  "
$ws.Range("E29").Value = "  +3.71%  "
$ws.Range("E30").Value = "  +5.10%  "
$ws.Range("E31").Value = "  +0.08%  "
$ws.Range("D32").Value = "'539.56"
$ws.Range("D32").Style = $ws.Range("B32").Style
$ws.Range("E32").Value = "  +6.95%  "
$ws.Range("E33").Value = "  +1.57%  "
$ws.Range("D34").Value = "'6.67"
$ws.Range("D34").Style = $ws.Range("B34").Style
$ws.Range("E34").Value = "  +5.59%  "
$ws.Range("E35").Value = "  -3.57%  "
$ws.Range("E36").Value = "  +2.12%  "
$ws.Range("D37").Value = "'20.85"
$ws.Range("D37").Style = $ws.Range("B37").Style
$ws.Range("E37").Value = "  +3.18%  "
$ws.Range("D38").Value = "'163.52"
$ws.Range("D38").Style = $ws.Range("B38").Style
$ws.Range("E38").Value = "  -0.17%  "
$ws.Range("E39").Value = "  -0.45%  "
$ws.Range("E40").Value = "  -0.01%  "
$ws.Range("D41").Value = "'171.23"
$ws.Range("D41").Style = $ws.Range("B41").Style
$ws.Range("E41").Value = "  +3.55%  "
$ws.Range("E42").Value = "  -0.02%  "
$ws.Range("D43").Value = "'42.49"
$ws.Range("D43").Style = $ws.Range("B43").Style
$ws.Range("E43").Value = "  +0.29%  "
$ws.Range("E44").Value = "  +2.53%  "
$ws.Range("D45").Value = "'0.0616"
$ws.Range("D45").Style = $ws.Range("B45").Style
$ws.Range("E45").Value = "  +0.72%  "
$ws.Range("E46").Value = "  +2.18%  "
$ws.Range("E47").Value = "  +4.41%  "
$ws.Range("E48").Value = "  +4.06%  "
$ws.Range("D49").Value = "'0.656"
$ws.Range("D49").Style = $ws.Range("B49").Style
$ws.Range("E49").Value = "  +1.51%  "
$ws.Range("D50").Value = "'20.95"
$ws.Range("D50").Style = $ws.Range("B50").Style
$ws.Range("E50").Value = "  +8.10%  "
$ws.Range("D51").Value = "'0.0992"
$ws.Range("D51").Style = $ws.Range("B51").Style
$ws.Range("E51").Value = "  +1.06%  "
